$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the DR= row label first so the new shared string ordering
# (DR=, Rserial, Discr) matches the authored workbook.
$ws.Range("D7").Value = "DR="

# New headers / labels for the discrete resistor calculation block
$ws.Range("D1").Value = "Rserial"
$ws.Range("E1").Value = 10000

$ws.Range("D2").Value = "Discr"
$ws.Range("E2").Value = 1023

# Update chosen temperature for R= row
$ws.Range("E3").Value = 160

# Add the discrete-value formula (single quotes so $E/$1 aren't expanded by PowerShell)
$ws.Range("E7").Formula = '=1023-1023/$E$1*E6'

# Keep selection on E3 as captured in the workbook view
$ws.Range("E3").Select()
